# Auto update Excel log
# Applies new sensor-log rows to ALERTS, PIR, Humidity, Proximity, and Camera sheets.

$wb = $excel.ActiveWorkbook

# --- ALERTS sheet (index 1) ---
$ws = $wb.Worksheets.Item(1)
$data = @(
    ,@("2026-01-30","18:23:06","18:00","Living Room","CRITICAL","FALL_DETECTED")
    ,@("2026-01-30","18:23:09","18:00","Living Room","CRITICAL","FALL_DETECTED")
)
$startRow = 6
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $rng = $ws.Range("A" + $r + ":F" + $r)
    $rng.NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# --- PIR sheet (index 2) ---
$ws = $wb.Worksheets.Item(2)
$data = @(
    ,@("2026-01-30","18:21:53","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:21:54","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:21:57","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:22:02","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:22:07","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:22:12","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:22:17","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:22:22","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:22:27","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:22:32","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:22:37","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:22:42","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:22:47","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:22:52","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:23:10","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:23:12","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:23:17","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:23:22","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:23:27","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:23:32","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:23:37","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:23:42","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:23:47","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:23:52","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:23:57","18:00","Bathroom","No Motion","Inactive")
    ,@("2026-01-30","18:24:02","18:00","Bathroom","No Motion","Inactive")
)
$startRow = 42
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $rng = $ws.Range("A" + $r + ":F" + $r)
    $rng.NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# --- Humidity sheet (index 3) ---
$ws = $wb.Worksheets.Item(3)
$data = @(
    ,@("2026-01-30","18:21:53","18:00","Bathroom","85.8%","Active")
    ,@("2026-01-30","18:21:54","18:00","Bathroom","86.7%","Active")
    ,@("2026-01-30","18:22:02","18:00","Bathroom","86.7%","Active")
    ,@("2026-01-30","18:22:12","18:00","Bathroom","86.6%","Active")
    ,@("2026-01-30","18:22:22","18:00","Bathroom","86.6%","Active")
    ,@("2026-01-30","18:22:32","18:00","Bathroom","86.7%","Active")
    ,@("2026-01-30","18:22:37","18:00","Bathroom","86.7%","Active")
    ,@("2026-01-30","18:22:42","18:00","Bathroom","86.7%","Active")
    ,@("2026-01-30","18:23:10","18:00","Bathroom","86.7%","Active")
    ,@("2026-01-30","18:23:12","18:00","Bathroom","86.6%","Active")
    ,@("2026-01-30","18:23:22","18:00","Bathroom","86.7%","Active")
    ,@("2026-01-30","18:23:27","18:00","Bathroom","86.7%","Active")
    ,@("2026-01-30","18:23:32","18:00","Bathroom","86.6%","Active")
    ,@("2026-01-30","18:23:43","18:00","Bathroom","86.7%","Active")
    ,@("2026-01-30","18:23:47","18:00","Bathroom","86.6%","Active")
    ,@("2026-01-30","18:23:53","18:00","Bathroom","86.6%","Active")
    ,@("2026-01-30","18:24:03","18:00","Bathroom","86.7%","Active")
)
$startRow = 32
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $rng = $ws.Range("A" + $r + ":F" + $r)
    $rng.NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# --- Proximity sheet (index 5) ---
$ws = $wb.Worksheets.Item(5)
$data = @(
    ,@("2026-01-30","18:23:19","18:00","Living Room Main Door","ENTER","User ENTERED Living Room Main Door")
    ,@("2026-01-30","18:23:21","18:00","Living Room Main Door","EXIT","User EXITED Living Room Main Door")
)
$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $rng = $ws.Range("A" + $r + ":F" + $r)
    $rng.NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# --- Camera sheet (index 7) ---
$ws = $wb.Worksheets.Item(7)
$data = @(
    ,@("2026-01-30","18:23:19","18:00","Living Room Main Door","Image Captured (ENTER)","Active")
    ,@("2026-01-30","18:23:21","18:00","Living Room Main Door","Image Captured (EXIT)","Active")
)
$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $rng = $ws.Range("A" + $r + ":F" + $r)
    $rng.NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
